$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L2").Value = 1.24
$ws.Range("N2").Value = 6.4
$ws.Range("S2").Value = 2.22
$ws.Range("U2").Value = 2.86
$ws.Range("F5").Value = 2.16
$ws.Range("G5").Value = 2.5
$ws.Range("I5").Value = 5
$ws.Range("J5").Value = 2.62
$ws.Range("K5").Value = 3.6
$ws.Range("L5").Value = 1.62
$ws.Range("M5").Value = 1.14
$ws.Range("N5").Value = 2.24
$ws.Range("O5").Value = 1.61
$ws.Range("P5").Value = 1.43
$ws.Range("R5").Value = 1.15
$ws.Range("S5").Value = 3.45
$ws.Range("T5").Value = 2.24
$ws.Range("U5").Value = 1.64
$ws.Range("V5").Value = 1.27
$ws.Range("W5").Value = 1.68
$ws.Range("F6").Value = 1.42
$ws.Range("G6").Value = 1.66
$ws.Range("H6").Value = 7
$ws.Range("K6").Value = 6.6
$ws.Range("N6").Value = 2.46
$ws.Range("W6").Value = 2.5
$ws.Range("U7").Value = 1.45
$ws.Range("F8").Value = 1.3
$ws.Range("G8").Value = 1.39
$ws.Range("H8").Value = 8
$ws.Range("J8").Value = 1.2
$ws.Range("N8").Value = 5.7
$ws.Range("W8").Value = 3.55
$ws.Range("F9").Value = 3.1
$ws.Range("G9").Value = 3.6
$ws.Range("K9").Value = 4
$ws.Range("N9").Value = 4
$ws.Range("O9").Value = 1.27
$ws.Range("P9").Value = 2.02
$ws.Range("Q9").Value = 1.78
$ws.Range("R9").Value = 1.41
$ws.Range("S9").Value = 3
$ws.Range("T9").Value = 1.67
$ws.Range("U9").Value = 2.22
$ws.Range("V9").Value = 1.69
$ws.Range("W9").Value = 1.38
$ws.Range("X9").Value = 20
$ws.Range("AC9").Value = 10
$ws.Range("AJ9").Value = 60
$ws.Range("AL9").Value = 980
$ws.Range("AM9").Value = 95
$ws.Range("P10").Value = 1.91
$ws.Range("Q10").Value = 1.48
$ws.Range("R10").Value = 1.35
$ws.Range("S10").Value = 2.58
$ws.Range("X10").Value = 26
$ws.Range("J11").Value = 3.9
$ws.Range("N12").Value = 2.98
$ws.Range("AC12").Value = 11
$ws.Range("F13").Value = 6.8
$ws.Range("J13").Value = 4.8
$ws.Range("L13").Value = 1.2
$ws.Range("R13").Value = 1.49
$ws.Range("S13").Value = 2.16
$ws.Range("H14").Value = 3.25
$ws.Range("W14").Value = 1.72
$ws.Range("F15").Value = 1.25
$ws.Range("H15").Value = 12
$ws.Range("I15").Value = 14
$ws.Range("J15").Value = 7.4
$ws.Range("K15").Value = 8
$ws.Range("N15").Value = 3.1
$ws.Range("P15").Value = 3.1
$ws.Range("Q15").Value = 1.37
$ws.Range("T15").Value = 1.65
$ws.Range("U15").Value = 1.75
$ws.Range("AF15").Value = 13
$ws.Range("AJ15").Value = 980
$ws.Range("AK15").Value = 980
$ws.Range("AL15").Value = 980
$ws.Range("AN15").Value = 4.6
$ws.Range("K16").Value = 4.2
$ws.Range("Q16").Value = 1.75
$ws.Range("T16").Value = 1.66
$ws.Range("Z16").Value = 980
$ws.Range("AF16").Value = 17
$ws.Range("AH16").Value = 21
$ws.Range("AJ16").Value = 980
$ws.Range("AL16").Value = 980
$ws.Range("AM16").Value = 100
$ws.Range("AN16").Value = 16
$ws.Range("F17").Value = 5.2
$ws.Range("I17").Value = 1.76
$ws.Range("J18").Value = 3.15
$ws.Range("L18").Value = 1.37
$ws.Range("O18").Value = 1.3
$ws.Range("Q18").Value = 2.02
$ws.Range("U18").Value = 2.04
$ws.Range("AA18").Value = 60
$ws.Range("AE18").Value = 44
$ws.Range("AK18").Value = 34
$ws.Range("AO18").Value = 42
$ws.Range("F21").Value = 5.1
$ws.Range("Y21").Value = 10.5
$ws.Range("AJ21").Value = 170
$ws.Range("AM21").Value = 140
$ws.Range("AN21").Value = 110
$ws.Range("F22").Value = 2.4
$ws.Range("Q22").Value = 2.48
$ws.Range("F23").Value = 2.1
$ws.Range("I23").Value = 3.75
$ws.Range("J23").Value = 3.55
$ws.Range("K23").Value = 5.8
$ws.Range("N23").Value = 2.28
$ws.Range("P23").Value = 2.26
$ws.Range("Q23").Value = 1.53
$ws.Range("V23").Value = 1.37
$ws.Range("H24").Value = 12
$ws.Range("I24").Value = 17
$ws.Range("J24").Value = 8
$ws.Range("L24").Value = 1.01
$ws.Range("N24").Value = 9.199999999999999
$ws.Range("O24").Value = 1.09
$ws.Range("P24").Value = 3.7
$ws.Range("Q24").Value = 1.3
$ws.Range("R24").Value = 2.12
$ws.Range("S24").Value = 1.75
$ws.Range("T24").Value = 1.73
$ws.Range("U24").Value = 2.12
$ws.Range("W24").Value = 5.1
$ws.Range("X24").Value = 65
$ws.Range("Z24").Value = 180
$ws.Range("AA24").Value = 550
$ws.Range("AC24").Value = 22
$ws.Range("AG24").Value = 13
$ws.Range("AH24").Value = 34
$ws.Range("AJ24").Value = 12
$ws.Range("AK24").Value = 13
$ws.Range("AN24").Value = 2.94
$ws.Range("AO24").Value = 1000
$ws.Range("L25").Value = 1.31
$ws.Range("O25").Value = 1.22
$ws.Range("T25").Value = 1.56
$ws.Range("U25").Value = 2.7
$ws.Range("Z25").Value = 22
